$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.890.16'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").Value = '2.600.60'
$ws.Range("E3").Value = '  -1.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.17'
$ws.Range("E5").Value = '  +3.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.50'
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  +4.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.81'
$ws.Range("E9").Value = '  +2.72%  '

$ws.Range("E10").Value = '  -0.99%  '

$ws.Range("E11").Value = '  +3.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.335'
$ws.Range("E12").Value = '  -0.59%  '

$ws.Range("D13").Value = '3.041.74'
$ws.Range("E13").Value = '  -2.15%  '

$ws.Range("D14").Value = '58.808.22'
$ws.Range("E14").Value = '  -0.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.89'
$ws.Range("E15").Value = '  -0.83%  '

$ws.Range("D16").Value = '2.602.68'
$ws.Range("E16").Value = '  -1.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  -1.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.46'
$ws.Range("E18").Value = '  +1.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '337.85'
$ws.Range("E19").Value = '  -0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.09'
$ws.Range("E20").Value = '  -2.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.14'
$ws.Range("E21").Value = '  -2.12%  '

$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.64'
$ws.Range("E23").Value = '  -0.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.427'
$ws.Range("E24").Value = '  +2.90%  '

$ws.Range("E25").Value = '  -0.32%  '

$ws.Range("E26").Value = '  -3.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.12'
$ws.Range("E27").Value = '  -2.05%  '

$ws.Range("D28").Value = '0.0₃0760'
$ws.Range("E28").Value = '  +2.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.67'
$ws.Range("E30").Value = '  +1.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.95'
$ws.Range("E31").Value = '  +1.99%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.03'
$ws.Range("E32").Value = '  +2.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.87'
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.93'
$ws.Range("E34").Value = '  -1.54%  '

$ws.Range("B35").Value = 'SuiNetwork'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.875'
$ws.Range("E35").Value = '  +4.79%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '37.22'
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.11'
$ws.Range("E37").Value = '  -1.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.45'
$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.827'
$ws.Range("E39").Value = '  -0.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.61'
$ws.Range("E40").Value = '  +0.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '285.65'
$ws.Range("E41").Value = '  -1.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.994'
$ws.Range("E42").Value = '  -0.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.594'
$ws.Range("E43").Value = '  -1.56%  '

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0952'
$ws.Range("E44").Value = '  +0.70%  '

$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.62'
$ws.Range("E45").Value = '  -0.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0532'
$ws.Range("E46").Value = '  -0.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0228'
$ws.Range("E47").Value = '  +0.97%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.62'
$ws.Range("E48").Value = '  +1.68%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.918.69'
$ws.Range("E49").Value = '  -2.60%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.99'
$ws.Range("E50").Value = '  -1.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '114.55'
$ws.Range("E51").Value = '  +3.20%  '
